# Apply the golden-test style updates described by the commit:
#   - Add a new "Abstract Title" paragraph style (based on Normal, next ->
#     Abstract) with centered, bold, colored 10pt text and spacing before.
#   - Tighten the "Abstract" style's space-before from 15pt (300) to 5pt (100).
#   - Give the ImportTok character style a bold green color.
#   - Give the BuiltInTok character style a green color.

$d = $word.ActiveDocument

# 1. New "Abstract Title" paragraph style.
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# 2. "Abstract" style: space-before 300 -> 100 (twips, i.e. 15pt -> 5pt).
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 3. "ImportTok" character style: bold, green (#008000).
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# 4. "BuiltInTok" character style: green (#008000).
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768

Write-Output "styles updated"
